# Apply "Add files via upload" changes:
# - add an "исполнители" (performers) column (E) with names
# - add a "role" column (F) with role/position for each performer
# - header cell E1 styled like the other headers
# - resize columns C (drop bestFit), E and F (new) to fit new content
# - update selection to H5
# - add a third stacked-bar chart series for "исполнители"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- New header cell (bold, like the other header cells) ----
$ws.Range("E1").Value = "исполнители"
$ws.Range("E1").Font.Bold = $true

# ---- New data: performers (E) and their roles (F) ----
# Values are entered in this exact order so new shared-string entries
# land at the same indices as in the source workbook.
$ws.Range("E2").Value = "Тхагазитов Эльдар"
$ws.Range("E5").Value = "Плотников Владимир"
$ws.Range("E3").Value = "Ерилеев Сергей"
$ws.Range("E4").Value = "Исмоилова Лейла"

$ws.Range("F2").Value = "разработчик "
$ws.Range("F5").Value = "проектировщик"
$ws.Range("F4").Value = "дизайнер"
$ws.Range("F3").Value = "разработчик "

# ---- Column widths ----
$ws.Columns("C").ColumnWidth = 16.42
$ws.Columns("E").ColumnWidth = 20.42
$ws.Columns("F").ColumnWidth = 14.75

# ---- Selection ----
$ws.Range("H5").Select()

# ---- Chart: add 3rd series "исполнители" (values from E2:E5) ----
$co = $ws.ChartObjects(1)
$chart = $co.Chart
$newSeries = $chart.SeriesCollection().NewSeries()
$newSeries.Format.Fill.Solid()
$newSeries.Format.Fill.ForeColor.RGB = 0xA5A5A5
$newSeries.Format.Line.Visible = $false
$newSeries.Formula = '=SERIES("исполнители",,Лист1!$E$2:$E$5,3)'
$chart.SeriesCollection(3).InvertIfNegative = $false
